$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zero-out the negative "Current Inventory" values (column F) on existing rows
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(5, 6).Value = 0

# Add a new test item as row 7
$ws.Cells.Item(7, 1).Value = "TEST01"
$ws.Cells.Item(7, 2).Value = "Test Item"
$ws.Cells.Item(7, 3).Value = "n/a"
$ws.Cells.Item(7, 4).Value = "TestVendor"
$ws.Cells.Item(7, 5).Value = 11.1
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 1

# The UPC text carries a trailing zero-width space, which would make the
# value look numeric if typed directly. Seed it with a non-numeric marker
# so Excel stores it as text, then swap the marker character for the
# real zero-width space via the Characters API (equivalent to editing the
# text in place), avoiding any forced "quote prefix" cell formatting.
$ws.Cells.Item(7, 8).Value = "123456789101X"
$ws.Cells.Item(7, 8).Characters(13, 1).Text = [char]0x200B
